$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) New title paragraph "S O P" inserted before the very first paragraph
#    ("BP or SP"), followed by a new blank paragraph.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphBefore()

$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.InsertAfter("S O P")
$titleRange.Font.Bold = $true
$titleRange.Font.BoldBi = $true
$titleRange.Font.Size = 36
$titleRange.Font.SizeBi = 36
$titlePara.Alignment = 1

$blankPara = $d.Paragraphs(2)
$blankPara.Range.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# 2) Fix the "Depp" typo -> "Deep" in the "1st Tank" line.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Depp", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Deep", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Fix the "TLT" typo -> "TILT" in the "No WAR" line.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("TLT", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TILT", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Append a new blank line, then "DBP - DSP Pattern 1-4", then another
#    blank line at the very end of the document.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()    # new blank paragraph (stays blank)

$blankBeforeDbp = $d.Paragraphs($d.Paragraphs.Count)
$blankBeforeDbp.Range.InsertParagraphAfter()   # new paragraph that will hold DBP text

$dbpPara = $d.Paragraphs($d.Paragraphs.Count)
$dbpPara.Range.InsertAfter("DBP " + [char]0x2013 + " DSP Pattern 1-4")
$dbpPara.Range.Font.Color = 0x47AD70
$dbpPara.Range.Font.Size = 24

$dbpPara.Range.InsertParagraphAfter()   # trailing blank paragraph at doc end

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host "$i -> [" $d.Paragraphs($i).Range.Text "]"
}
